# Master Data Tables - Test Data / master-template_type.xlsx
# "Updated Master data as per 16th May Refresh"
#
# 1) Rename the three "otp-*" template codes to the new "ida-auth-otp-*-template"
#    codes (each code appears 3x - once per language row: eng, ara, fra).
# 2) Append 12 new rows (125-136) for the new "consent" and
#    "auth-otp-*-template" master data entries (eng/ara/fra each).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename otp-* codes -> ida-auth-otp-*-template --------------------
# English block (rows 5-7)
$ws.Range("A5").Value  = "ida-auth-otp-email-content-template"
$ws.Range("A6").Value  = "ida-auth-otp-email-subject-template"
$ws.Range("A7").Value  = "ida-auth-otp-sms-template"

# Arabic block (rows 11-13)
$ws.Range("A11").Value = "ida-auth-otp-email-content-template"
$ws.Range("A12").Value = "ida-auth-otp-email-subject-template"
$ws.Range("A13").Value = "ida-auth-otp-sms-template"

# French block (rows 17-19)
$ws.Range("A17").Value = "ida-auth-otp-email-content-template"
$ws.Range("A18").Value = "ida-auth-otp-email-subject-template"
$ws.Range("A19").Value = "ida-auth-otp-sms-template"

# --- 2) Append new rows 125-136 -------------------------------------------
$newRows = @(
    @("consent",                          "Consent",                           "eng", $true, "superadmin", "now()"),
    @("consent",                          "موافقة",                            "ara", $true, "superadmin", "now()"),
    @("consent",                          "Consentement",                      "fra", $true, "superadmin", "now()"),
    @("auth-otp-email-subject-template",  "Auth OTP Email Subject Template",   "eng", $true, "superadmin", "now()"),
    @("auth-otp-email-subject-template",  "مصادقة OTP قالب موضوع",              "ara", $true, "superadmin", "now()"),
    @("auth-otp-email-subject-template",  "Modèle dobjet de-mail Auth OTP",    "fra", $true, "superadmin", "now()"),
    @("auth-otp-email-content-template",  "Auth OTP Email Content Template",   "eng", $true, "superadmin", "now()"),
    @("auth-otp-email-content-template",  "مصادقة OTP قالب محتوى",              "ara", $true, "superadmin", "now()"),
    @("auth-otp-email-content-template",  "Auth OTP Email ContentTemplate",    "fra", $true, "superadmin", "now()"),
    @("auth-otp-sms-template",            "Auth OTP SMS Template",              "eng", $true, "superadmin", "now()"),
    @("auth-otp-sms-template",            "مصادقة قالب رسالة OTP",              "ara", $true, "superadmin", "now()"),
    @("auth-otp-sms-template",            "Modèle SMS OTP Auth",               "fra", $true, "superadmin", "now()")
)

$startRow = 125
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
}

# --- Update the sheet selection to mirror the post-edit Excel state -------
# (the author's Excel session left the cursor on the row right after the
# last data row, with the rest of the column selected)
$ws.Range("A137:XFD1048576").Select()
